$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C18) from 45208 to 45212
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}

# Update the hyperlink formulas in row 2 (S2:Y2) with new file names
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/artfynd/A 30779-2023 artfynd.xlsx", "A 30779-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/kartor/A 30779-2023 karta.png", "A 30779-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/knärot/A 30779-2023 karta knärot.png", "A 30779-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/klagomål/A 30779-2023 fsc-klagomål.docx", "A 30779-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/klagomålsmail/A 30779-2023 fsc-klagomål mail.docx", "A 30779-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/tillsyn/A 30779-2023 tillsynsbegäran.docx", "A 30779-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0883/ti,llsynsmail/A 30779-2023 tillsynsbegäran mail.docx", "A 30779-2023")'
